$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.99999999273189988,
    0.9927553246069557,
    0.97107650281558056,
    0.96095820562386169,
    0.95124068021528729,
    0.92734985993311303,
    0.92575228307933943,
    0.92422099241327693,
    0.9255339419277766,
    0.92782001582896734,
    0.92829634835634445,
    0.9255165682466886,
    0.91422965683109236,
    0.91006237024434655,
    0.90747092679941765,
    0.90496446448529122,
    0.90125659055898577,
    0.90014768995765149,
    0.98896659261537101,
    0.98184954099625865,
    0.98045106767183832,
    0.97918657028452194,
    0.9640407785155608,
    0.95101895442310425,
    0.94456180672465873,
    0.92369274596818196,
    0.91884409252079158,
    0.89736435983527874,
    0.88209052291489276,
    0.87551963892766438,
    0.86786554428256046,
    0.8661862302986949,
    0.86566623293025624
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
